$d = $word.ActiveDocument

$d.Content.Find.Execute("Medinis korpusas", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Žaliavų stoka", 2)

$d.Content.Find.Execute("Irdamas medis gamtoje neišskleidžia nuodingų medžiagų.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Produkto korpusui taip pat galima naudoti ir medį, kurio stokos nėra. Taip pat naudoti perdirbtą plastiką.", 2)
